$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply refreshed cryptocurrency price / 1h-volume data, and the Frax<->Aave row swap (rows 48-49).
# D-column price strings are forced to Text (NumberFormat "@") before assignment so numeric-looking
# values like "1.003" or "214.89" stay literal text (matching the source feed formatting) instead of
# being coerced into floating point numbers by Excel; Style is reset to Normal afterwards so no stray
# cell formatting is introduced.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '26.011.59'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.34%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.641.75'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -0.09%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '214.89'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.24%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.5085'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.71%  '
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2561'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.29%  '
$ws.Range('E9').Value = '  -0.76%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.55'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.43%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07776'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.16%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '4.285'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.14%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.649.89'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.57%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.5423'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.13%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '64.09'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -1.13%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0₅7698'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -2.04%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '26.036.59'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.29%  '
$ws.Range('E18').Value = '  -0.30%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '198.86'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.26%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.422'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.68%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '9.897'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.80%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.043'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.05%  '
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.869'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '141.10'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.84%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.1191'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +4.18%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '6.810'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.62%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '15.65'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.41%  '
$ws.Range('E29').Value = '  -0.72%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.04890'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.64%  '
$ws.Range('E31').Value = '  -0.45%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.162'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -1.13%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.525'
$ws.Range('D33').Style = "Normal"
$ws.Range('E34').Value = '  -0.37%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.9029'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.98%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.583'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.95%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.142.93'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.5450'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.82%  '
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.002'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.37%  '
$ws.Range('E41').Value = '  -1.22%  '
$ws.Range('E42').Value = '  +7.60%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.8101'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.93%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '99.37'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.09%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '5.399'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -5.15%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.4534'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.18%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '54.98'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.72%  '
$ws.Range('B49').Value = 'Frax'
$ws.Range('C49').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.001'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.50%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.05090'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.26%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.004'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.28%  '
